$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.312.99"
Set-TextValue $ws.Range("E2") "  +0.28%  "
Set-TextValue $ws.Range("D3") "1.869.84"
Set-TextValue $ws.Range("E3") "  -0.32%  "
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "235.60"
Set-TextValue $ws.Range("E5") "  +0.26%  "
Set-TextValue $ws.Range("E6") "  -0.06%  "
Set-TextValue $ws.Range("D7") "0.4684"
Set-TextValue $ws.Range("E7") "  +0.22%  "
Set-TextValue $ws.Range("E8") "  +0.83%  "
Set-TextValue $ws.Range("D9") "0.06538"
Set-TextValue $ws.Range("E9") "  -0.38%  "
Set-TextValue $ws.Range("D10") "21.35"
Set-TextValue $ws.Range("E10") "  +4.92%  "
Set-TextValue $ws.Range("D11") "0.07868"
Set-TextValue $ws.Range("E11") "  +1.49%  "
Set-TextValue $ws.Range("E12") "  +0.50%  "
Set-TextValue $ws.Range("D13") "1.862.46"
Set-TextValue $ws.Range("E13") "  -0.75%  "
Set-TextValue $ws.Range("D14") "5.102"
Set-TextValue $ws.Range("E14") "  +0.80%  "
Set-TextValue $ws.Range("D15") "0.6759"
Set-TextValue $ws.Range("E15") "  +0.75%  "
Set-TextValue $ws.Range("D16") "279.34"
Set-TextValue $ws.Range("D17") "30.309.53"
Set-TextValue $ws.Range("E17") "  +0.19%  "
Set-TextValue $ws.Range("D18") "1.000"
Set-TextValue $ws.Range("E18") "  +0.00%  "
Set-TextValue $ws.Range("D19") "5.506"
Set-TextValue $ws.Range("E19") "  +2.55%  "
Set-TextValue $ws.Range("D20") "12.72"
Set-TextValue $ws.Range("E20") "  +1.17%  "
Set-TextValue $ws.Range("D21") "2.108.33"
Set-TextValue $ws.Range("E21") "  -0.88%  "
Set-TextValue $ws.Range("D22") "0.000007303"
Set-TextValue $ws.Range("E22") "  +1.00%  "
Set-TextValue $ws.Range("D23") "1.000"
Set-TextValue $ws.Range("E23") "  -0.16%  "
Set-TextValue $ws.Range("D24") "6.166"
Set-TextValue $ws.Range("E24") "  +0.14%  "
Set-TextValue $ws.Range("D25") "165.38"
Set-TextValue $ws.Range("E25") "  -1.43%  "
Set-TextValue $ws.Range("D26") "9.173"
Set-TextValue $ws.Range("E26") "  -1.51%  "
Set-TextValue $ws.Range("E27") "  +0.12%  "
Set-TextValue $ws.Range("D28") "1.932"
Set-TextValue $ws.Range("E28") "  -2.19%  "
Set-TextValue $ws.Range("E29") "  +0.00%  "
Set-TextValue $ws.Range("D30") "0.09637"
Set-TextValue $ws.Range("E30") "  +0.05%  "
Set-TextValue $ws.Range("D31") "4.372"
Set-TextValue $ws.Range("E31") "  +0.18%  "
Set-TextValue $ws.Range("E32") "  +0.65%  "
Set-TextValue $ws.Range("D33") "4.097"
Set-TextValue $ws.Range("E33") "  -0.15%  "
Set-TextValue $ws.Range("E34") "  +1.08%  "
Set-TextValue $ws.Range("D35") "1.129"
Set-TextValue $ws.Range("E35") "  +3.60%  "
Set-TextValue $ws.Range("D36") "0.7065"
Set-TextValue $ws.Range("E36") "  +0.92%  "
Set-TextValue $ws.Range("D37") "2.722"
Set-TextValue $ws.Range("E37") "  +0.17%  "
Set-TextValue $ws.Range("D38") "0.01856"
Set-TextValue $ws.Range("E38") "  -0.29%  "
Set-TextValue $ws.Range("D39") "6.278"
Set-TextValue $ws.Range("E39") "  -3.83%  "
Set-TextValue $ws.Range("D40") "2.529"
Set-TextValue $ws.Range("E40") "  +0.20%  "
Set-TextValue $ws.Range("D41") "74.02"
Set-TextValue $ws.Range("E41") "  +3.13%  "
Set-TextValue $ws.Range("E42") "  +0.03%  "
Set-TextValue $ws.Range("D43") "0.8493"
Set-TextValue $ws.Range("E43") "  -1.37%  "
Set-TextValue $ws.Range("D44") "0.4182"
Set-TextValue $ws.Range("E44") "  +0.32%  "
Set-TextValue $ws.Range("D47") "7.170"
Set-TextValue $ws.Range("E47") "  -0.31%  "
Set-TextValue $ws.Range("D48") "9.252"
Set-TextValue $ws.Range("E48") "  +1.19%  "
Set-TextValue $ws.Range("D49") "936.00"
Set-TextValue $ws.Range("E49") "  -4.79%  "
Set-TextValue $ws.Range("E50") "  +1.19%  "
Set-TextValue $ws.Range("E51") "  -1.63%  "

# Row 45/46 swap (Quant <-> PaxDollar)
Set-TextValue $ws.Range("B45") "PaxDollar"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D45") "1.000"
Set-TextValue $ws.Range("E45") "  -0.02%  "

Set-TextValue $ws.Range("B46") "Quant"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D46") "103.93"
Set-TextValue $ws.Range("E46") "  +1.02%  "
